$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

# --- Add new row 3 data -----------------------------------------------
$ws.Range("A3").Value = "TC_LoginTest_02"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "sammca87@gmail.com"
$ws.Range("D3").Value = "asdfasdf"

# --- Stash the existing "Hyperlink" look-and-feel from row 2 so we can
#     re-apply it after Hyperlinks.Add() (which stamps its own style). ---
$ws.Range("C2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("F3").PasteSpecial(-4122)

# --- Rebuild the row-2 hyperlinks so rId1/rId2 swap targets (D2 first,
#     then C2), matching the re-saved relationship order. -------------
$null = $ws.Range("C2").Hyperlinks.Delete()
$null = $ws.Range("D2").Hyperlinks.Delete()

$null = $ws.Hyperlinks.Add($ws.Range("D2"), "mailto:yuvi@666")
$null = $ws.Hyperlinks.Add($ws.Range("C2"), "mailto:sammca87@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("D3"), "mailto:yuvi@666", "", "", "yuvi@666")
$null = $ws.Hyperlinks.Add($ws.Range("C3"), "mailto:sammca87@gmail.com")

# Hyperlinks.Add() overwrites the cell text with the link target / the
# supplied display text -- restore the real cell values.
$ws.Range("C2").Value = "sammca87@gmail.com"
$ws.Range("D2").Value = "yuvi@666"
$ws.Range("C3").Value = "sammca87@gmail.com"
$ws.Range("D3").Value = "asdfasdf"

# Re-apply the original Hyperlink formatting (font/underline/colour) so
# the cells keep using the workbook's existing style instead of a new one.
$ws.Range("F2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# Clean up the helper cells used to stash formatting.
$ws.Range("F2:F3").Clear()

# --- Selection follows the newly-entered cell -------------------------
$null = $ws.Range("D3").Select()
